$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 675
$ws.Cells.Item(4, 6).Value = 1987
$ws.Cells.Item(5, 6).Value = 5782
$ws.Cells.Item(6, 6).Value = 1624
$ws.Cells.Item(7, 6).Value = 171
$ws.Cells.Item(8, 6).Value = 3270
$ws.Cells.Item(10, 6).Value = 46
$ws.Cells.Item(11, 6).Value = 1372
$ws.Cells.Item(12, 6).Value = 4564
$ws.Cells.Item(13, 6).Value = 1088
$ws.Cells.Item(14, 6).Value = 1719
$ws.Cells.Item(18, 6).Value = 53
$ws.Cells.Item(19, 6).Value = 181
$ws.Cells.Item(20, 6).Value = 156
$ws.Cells.Item(21, 6).Value = 1030
$ws.Cells.Item(22, 6).Value = 308
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(27, 6).Value = 213
$ws.Cells.Item(29, 6).Value = 1124
$ws.Cells.Item(30, 6).Value = 412
$ws.Cells.Item(31, 6).Value = 89
$ws.Cells.Item(33, 6).Value = 392
$ws.Cells.Item(34, 6).Value = 965
$ws.Cells.Item(35, 6).Value = 17
$ws.Cells.Item(36, 6).Value = 1754
$ws.Cells.Item(37, 6).Value = 2257
$ws.Cells.Item(38, 6).Value = 1053
$ws.Cells.Item(40, 6).Value = 278
$ws.Cells.Item(41, 6).Value = 639
$ws.Cells.Item(42, 6).Value = 379
$ws.Cells.Item(43, 6).Value = 40
$ws.Cells.Item(44, 6).Value = 677
$ws.Cells.Item(45, 6).Value = 32
$ws.Cells.Item(46, 6).Value = 446
$ws.Cells.Item(47, 6).Value = 413
$ws.Cells.Item(48, 6).Value = 232

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(11, 6).Value = 157
$ws.Cells.Item(11, 7).Value = 252
$ws.Cells.Item(22, 6).Value = 4

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 783

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 783
$ws.Cells.Item(3, 6).Value = 675
$ws.Cells.Item(5, 6).Value = 1987
$ws.Cells.Item(6, 6).Value = 5782
$ws.Cells.Item(7, 6).Value = 1624
$ws.Cells.Item(8, 6).Value = 171
$ws.Cells.Item(9, 6).Value = 3270
$ws.Cells.Item(10, 6).Value = 46
$ws.Cells.Item(11, 6).Value = 1372
$ws.Cells.Item(12, 6).Value = 4564
$ws.Cells.Item(13, 6).Value = 1088
$ws.Cells.Item(14, 6).Value = 1719
$ws.Cells.Item(19, 6).Value = 53
$ws.Cells.Item(20, 6).Value = 181
$ws.Cells.Item(21, 6).Value = 156
$ws.Cells.Item(22, 6).Value = 157
$ws.Cells.Item(22, 7).Value = 252
$ws.Cells.Item(23, 6).Value = 1030
$ws.Cells.Item(24, 6).Value = 308
$ws.Cells.Item(27, 6).Value = 213
$ws.Cells.Item(29, 6).Value = 1124
$ws.Cells.Item(30, 6).Value = 412
$ws.Cells.Item(31, 6).Value = 89
$ws.Cells.Item(33, 6).Value = 965
$ws.Cells.Item(34, 6).Value = 1754
$ws.Cells.Item(35, 6).Value = 2257
$ws.Cells.Item(36, 6).Value = 1053
$ws.Cells.Item(40, 6).Value = 278
$ws.Cells.Item(41, 6).Value = 639
$ws.Cells.Item(42, 6).Value = 379
$ws.Cells.Item(43, 6).Value = 677
$ws.Cells.Item(44, 6).Value = 446
$ws.Cells.Item(45, 6).Value = 413
$ws.Cells.Item(46, 6).Value = 232
$ws.Cells.Item(47, 6).Value = 4
